$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibition) ---
$ws = $wb.Worksheets.Item(1)
$ws1 = $ws
# Row 2 (pre-insert update)
$ws1.Range("A2").Value = 1
$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = '2024-10-26'
$ws1.Range("C2").Value = '合肥·W·A第五人格同人only2.0'
$ws1.Range("D2").Value = '莲花路与石门路交口西北角（尚泽大都会B座四楼） 格律诗婚礼艺术中心(经开店)'
$ws1.Range("E2").Value = '2024.10.26 09:30-10.26 17:00'
$ws1.Range("F2").Value = 454
$ws1.Range("G2").Value = '已售罄'
$ws1.Range("H2").Value = 'https://show.bilibili.com/platform/detail.html?id=91123'
$ws1.Range("I2").Value = '//i2.hdslb.com/bfs/openplatform/202408/YqXHTFM81724066565119.png'
$ws1.Rows.Item(3).Insert()
# Row 3
$ws1.Range("A3").Value = 2
$ws1.Range("B3").NumberFormat = "@"
$ws1.Range("B3").Value = '2024-11-02'
$ws1.Range("C3").Value = '合肥·之心城购物中心-2024漫趣地带嘉年华（免费）'
$ws1.Range("D3").Value = '长江西路189号 之心城'
$ws1.Range("E3").Value = '2024.11.02 10:00-11.03 22:00'
$ws1.Range("F3").Value = 1
$ws1.Range("G3").Value = 30
$ws1.Range("H3").Value = 'https://show.bilibili.com/platform/detail.html?id=93887'
$ws1.Range("I3").Value = '//i1.hdslb.com/bfs/openplatform/202410/JbPnmkCQ1729673353796.jpeg'
# Row 4
$ws1.Range("A4").Value = 3
$ws1.Range("B4").NumberFormat = "@"
$ws1.Range("B4").Value = '2024-11-09'
$ws1.Range("C4").Value = '安徽·崩坏同人only 爱莉希雅同人生日会'
$ws1.Range("D4").Value = '徽州大道与扬子江路口天琅百老汇一楼123号 禧棠捌号XITANGBH禧宴中心'
$ws1.Range("E4").Value = '2024.11.09 12:00-11.09 22:00'
$ws1.Range("F4").Value = 35
$ws1.Range("G4").Value = '不可售'
$ws1.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=93461'
$ws1.Range("I4").Value = '//i2.hdslb.com/bfs/openplatform/202410/VnEQZYTQ1728892328769.png'
# Row 5
$ws1.Range("A5").Value = 4
$ws1.Range("B5").NumberFormat = "@"
$ws1.Range("B5").Value = '2024-11-16'
$ws1.Range("C5").Value = '合肥·11.16合肥耽美同人only'
$ws1.Range("D5").Value = '金寨路与天堂窄路交叉口 梵木艺术中心'
$ws1.Range("E5").Value = '2024.11.16 10:00-11.16 17:00'
$ws1.Range("F5").Value = 75
$ws1.Range("G5").Value = 68
$ws1.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=93612'
$ws1.Range("I5").Value = '//i1.hdslb.com/bfs/openplatform/202410/zRIIv4H81729147685895.jpeg'
# Row 6
$ws1.Range("A6").Value = 5
$ws1.Range("B6").NumberFormat = "@"
$ws1.Range("B6").Value = '2024-11-16'
$ws1.Range("C6").Value = '合肥·第九届环形宇宙动漫游戏嘉年华'
$ws1.Range("D6").Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$ws1.Range("E6").Value = '2024.11.16 09:30-11.17 17:00'
$ws1.Range("F6").Value = 5212
$ws1.Range("G6").Value = 72
$ws1.Range("H6").Value = 'https://show.bilibili.com/platform/detail.html?id=92565'
$ws1.Range("I6").Value = '//i1.hdslb.com/bfs/openplatform/202410/sxfiGFBq1728715876124.jpeg'
# Row 7
$ws1.Range("A7").Value = 6
$ws1.Range("B7").NumberFormat = "@"
$ws1.Range("B7").Value = '2024-11-17'
$ws1.Range("C7").Value = '合肥·MAX特摄同人only2.0'
$ws1.Range("D7").Value = '桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间'
$ws1.Range("E7").Value = '2024.11.17 10:00-11.17 18:00'
$ws1.Range("F7").Value = 176
$ws1.Range("G7").Value = 60
$ws1.Range("H7").Value = 'https://show.bilibili.com/platform/detail.html?id=92498'
$ws1.Range("I7").Value = '//i1.hdslb.com/bfs/openplatform/202410/ccEfc1521728888008037.jpeg'
# Row 8
$ws1.Range("A8").Value = 7
$ws1.Range("B8").NumberFormat = "@"
$ws1.Range("B8").Value = '2024-11-23'
$ws1.Range("C8").Value = '合肥·九号幻想动漫游戏嘉年华'
$ws1.Range("D8").Value = '龙岗路与淮南路交口东北角 合肥市青少年活动中心'
$ws1.Range("E8").Value = '2024.11.23 09:00-11.24 17:30'
$ws1.Range("F8").Value = 78
$ws1.Range("G8").Value = 69
$ws1.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=93609'
$ws1.Range("I8").Value = '//i1.hdslb.com/bfs/openplatform/202410/MLzZJXDx1729501364141.jpeg'
# Row 9
$ws1.Range("A9").Value = 8
$ws1.Range("B9").NumberFormat = "@"
$ws1.Range("B9").Value = '2024-11-30'
$ws1.Range("C9").Value = '合肥·风月引代号鸢同人only'
$ws1.Range("D9").Value = '徽州大道与杨子江路交口天琅百老汇一楼123号 禧棠捌号禧宴中心（滨湖店）'
$ws1.Range("E9").Value = '2024.11.30 10:00-11.30 21:00'
$ws1.Range("F9").Value = 98
$ws1.Range("G9").Value = 55
$ws1.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=93322'
$ws1.Range("I9").Value = '//i1.hdslb.com/bfs/openplatform/202409/Tu5YLbGx1727179854562.jpeg'
# Row 10
$ws1.Range("A10").Value = 9
$ws1.Range("B10").NumberFormat = "@"
$ws1.Range("B10").Value = '2024-12-07'
$ws1.Range("C10").Value = '合肥·心动恋章·冬日序国乙&代号鸢同人only'
$ws1.Range("D10").Value = '上海路与迎淮路交口向东200米(云峯中心一楼) 费加罗宴会艺术中心(省府店)'
$ws1.Range("E10").Value = '2024.12.07 12:00-12.07 21:00'
$ws1.Range("F10").Value = 344
$ws1.Range("G10").Value = 77
$ws1.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=93319'
$ws1.Range("I10").Value = '//i0.hdslb.com/bfs/openplatform/202409/KtMLL8ZO1727684987784.jpeg'
# Row 11
$ws1.Range("A11").Value = 10
$ws1.Range("B11").NumberFormat = "@"
$ws1.Range("B11").Value = '2024-12-08'
$ws1.Range("C11").Value = '合肥·星光国潮动漫游戏嘉年华'
$ws1.Range("D11").Value = '北二环与新蚌埠路交汇处 蓝金湾大酒店'
$ws1.Range("E11").Value = '2024.12.08 10:00-12.08 17:00'
$ws1.Range("F11").Value = 8
$ws1.Range("G11").Value = 39.9
$ws1.Range("H11").Value = 'https://show.bilibili.com/platform/detail.html?id=93801'
$ws1.Range("I11").Value = '//i0.hdslb.com/bfs/openplatform/202410/ubX6VZ841729253636894.png'
# Row 12
$ws1.Range("A12").Value = 11
$ws1.Range("B12").NumberFormat = "@"
$ws1.Range("B12").Value = '2025-02-03'
$ws1.Range("C12").Value = '合肥·皖萌次元青年文化节'
$ws1.Range("D12").Value = '凤淮路与公园路交叉口南行50米路西 庐阳区全民健身中心'
$ws1.Range("E12").Value = '2025.02.03 10:00-02.04 17:30'
$ws1.Range("F12").Value = 63
$ws1.Range("G12").Value = 39.9
$ws1.Range("H12").Value = 'https://show.bilibili.com/platform/detail.html?id=93589'
$ws1.Range("I12").Value = '//i0.hdslb.com/bfs/openplatform/202410/GjWiXfOf1729133962063.jpeg'

# --- Sheet 4: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item(4)
$ws4 = $ws
# Row 2 (pre-insert update)
$ws4.Range("A2").Value = 1
$ws4.Range("B2").NumberFormat = "@"
$ws4.Range("B2").Value = '2024-10-26'
$ws4.Range("C2").Value = '合肥·W·A第五人格同人only2.0'
$ws4.Range("D2").Value = '莲花路与石门路交口西北角（尚泽大都会B座四楼） 格律诗婚礼艺术中心(经开店)'
$ws4.Range("E2").Value = '2024.10.26 09:30-10.26 17:00'
$ws4.Range("F2").Value = 454
$ws4.Range("G2").Value = '已售罄'
$ws4.Range("H2").Value = 'https://show.bilibili.com/platform/detail.html?id=91123'
$ws4.Range("I2").Value = '//i2.hdslb.com/bfs/openplatform/202408/YqXHTFM81724066565119.png'
# Row 3 (pre-insert update)
$ws4.Range("A3").Value = 2
$ws4.Range("B3").NumberFormat = "@"
$ws4.Range("B3").Value = '2024-10-26'
$ws4.Range("C3").Value = '合肥·《四月是你的谎言》—“公生”与“薰”的钢琴小提琴唯美经典音乐集'
$ws4.Range("D3").Value = '徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院'
$ws4.Range("E3").Value = '2024.10.26 19:30-10.26 21:00'
$ws4.Range("F3").Value = 127
$ws4.Range("G3").Value = '不可售'
$ws4.Range("H3").Value = 'https://show.bilibili.com/platform/detail.html?id=90322'
$ws4.Range("I3").Value = '//i2.hdslb.com/bfs/openplatform/202408/BiVgXUKH1722824304648.jpeg'
$ws4.Rows.Item(4).Insert()
# Row 4
$ws4.Range("A4").Value = 3
$ws4.Range("B4").NumberFormat = "@"
$ws4.Range("B4").Value = '2024-11-02'
$ws4.Range("C4").Value = '合肥·之心城购物中心-2024漫趣地带嘉年华（免费）'
$ws4.Range("D4").Value = '长江西路189号 之心城'
$ws4.Range("E4").Value = '2024.11.02 10:00-11.03 22:00'
$ws4.Range("F4").Value = 1
$ws4.Range("G4").Value = 30
$ws4.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=93887'
$ws4.Range("I4").Value = '//i1.hdslb.com/bfs/openplatform/202410/JbPnmkCQ1729673353796.jpeg'
# Row 5
$ws4.Range("A5").Value = 4
$ws4.Range("B5").NumberFormat = "@"
$ws4.Range("B5").Value = '2024-11-08'
$ws4.Range("C5").Value = '合肥·松井祐贵 2024《阳光之旅》指弹吉他音乐会'
$ws4.Range("D5").Value = '宁国南路与水阳江路交口罍街二期15号楼安徽原创音乐基地3楼 OTW LIVEHOUSE'
$ws4.Range("E5").Value = '2024.11.08 19:30-11.08 21:00'
$ws4.Range("F5").Value = 2
$ws4.Range("G5").Value = 220
$ws4.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=92768'
$ws4.Range("I5").Value = '//i1.hdslb.com/bfs/openplatform/202409/OU2qWxgM1727082424391.jpeg'
# Row 6
$ws4.Range("A6").Value = 5
$ws4.Range("B6").NumberFormat = "@"
$ws4.Range("B6").Value = '2024-11-08'
$ws4.Range("C6").Value = '合肥·豫章D乐团-《蓬莱乐，万物生》——传统×先锋 疗愈音乐会'
$ws4.Range("D6").Value = '徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院'
$ws4.Range("E6").Value = '2024.11.08 19:30-11.08 21:00'
$ws4.Range("F6").Value = 1
$ws4.Range("G6").Value = 79.90000000000001
$ws4.Range("H6").Value = 'https://show.bilibili.com/platform/detail.html?id=92957'
$ws4.Range("I6").Value = '//i0.hdslb.com/bfs/openplatform/202409/uifvAByr1727253170481.jpeg'
# Row 7
$ws4.Range("A7").Value = 6
$ws4.Range("B7").NumberFormat = "@"
$ws4.Range("B7").Value = '2024-11-09'
$ws4.Range("C7").Value = '合肥·一生必听的钢琴曲—“从巴赫 · 莫扎特到肖邦 · 李斯特”钢琴圣手谭小棠独奏音乐会'
$ws4.Range("D7").Value = '徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院'
$ws4.Range("E7").Value = '2024.11.09 19:30-11.09 21:00'
$ws4.Range("F7").Value = 7
$ws4.Range("G7").Value = 80
$ws4.Range("H7").Value = 'https://show.bilibili.com/platform/detail.html?id=90593'
$ws4.Range("I7").Value = '//i2.hdslb.com/bfs/openplatform/202408/SYfLxnO21723442234232.jpeg'
# Row 8
$ws4.Range("A8").Value = 7
$ws4.Range("B8").NumberFormat = "@"
$ws4.Range("B8").Value = '2024-11-09'
$ws4.Range("C8").Value = '安徽·崩坏同人only 爱莉希雅同人生日会'
$ws4.Range("D8").Value = '徽州大道与扬子江路口天琅百老汇一楼123号 禧棠捌号XITANGBH禧宴中心'
$ws4.Range("E8").Value = '2024.11.09 12:00-11.09 22:00'
$ws4.Range("F8").Value = 35
$ws4.Range("G8").Value = '不可售'
$ws4.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=93461'
$ws4.Range("I8").Value = '//i2.hdslb.com/bfs/openplatform/202410/VnEQZYTQ1728892328769.png'
# Row 9
$ws4.Range("A9").Value = 8
$ws4.Range("B9").NumberFormat = "@"
$ws4.Range("B9").Value = '2024-11-16'
$ws4.Range("C9").Value = '合肥·11.16合肥耽美同人only'
$ws4.Range("D9").Value = '金寨路与天堂窄路交叉口 梵木艺术中心'
$ws4.Range("E9").Value = '2024.11.16 10:00-11.16 17:00'
$ws4.Range("F9").Value = 75
$ws4.Range("G9").Value = 68
$ws4.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=93612'
$ws4.Range("I9").Value = '//i1.hdslb.com/bfs/openplatform/202410/zRIIv4H81729147685895.jpeg'
# Row 10
$ws4.Range("A10").Value = 9
$ws4.Range("B10").NumberFormat = "@"
$ws4.Range("B10").Value = '2024-11-16'
$ws4.Range("C10").Value = '合肥·第九届环形宇宙动漫游戏嘉年华'
$ws4.Range("D10").Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$ws4.Range("E10").Value = '2024.11.16 09:30-11.17 17:00'
$ws4.Range("F10").Value = 5212
$ws4.Range("G10").Value = 72
$ws4.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=92565'
$ws4.Range("I10").Value = '//i1.hdslb.com/bfs/openplatform/202410/sxfiGFBq1728715876124.jpeg'
# Row 11
$ws4.Range("A11").Value = 10
$ws4.Range("B11").NumberFormat = "@"
$ws4.Range("B11").Value = '2024-11-17'
$ws4.Range("C11").Value = '合肥·MAX特摄同人only2.0'
$ws4.Range("D11").Value = '桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间'
$ws4.Range("E11").Value = '2024.11.17 10:00-11.17 18:00'
$ws4.Range("F11").Value = 176
$ws4.Range("G11").Value = 60
$ws4.Range("H11").Value = 'https://show.bilibili.com/platform/detail.html?id=92498'
$ws4.Range("I11").Value = '//i1.hdslb.com/bfs/openplatform/202410/ccEfc1521728888008037.jpeg'
# Row 12
$ws4.Range("A12").Value = 11
$ws4.Range("B12").NumberFormat = "@"
$ws4.Range("B12").Value = '2024-11-23'
$ws4.Range("C12").Value = '合肥·九号幻想动漫游戏嘉年华'
$ws4.Range("D12").Value = '龙岗路与淮南路交口东北角 合肥市青少年活动中心'
$ws4.Range("E12").Value = '2024.11.23 09:00-11.24 17:30'
$ws4.Range("F12").Value = 78
$ws4.Range("G12").Value = 69
$ws4.Range("H12").Value = 'https://show.bilibili.com/platform/detail.html?id=93609'
$ws4.Range("I12").Value = '//i1.hdslb.com/bfs/openplatform/202410/MLzZJXDx1729501364141.jpeg'
# Row 13
$ws4.Range("A13").Value = 12
$ws4.Range("B13").NumberFormat = "@"
$ws4.Range("B13").Value = '2024-11-30'
$ws4.Range("C13").Value = '合肥·风月引代号鸢同人only'
$ws4.Range("D13").Value = '徽州大道与杨子江路交口天琅百老汇一楼123号 禧棠捌号禧宴中心（滨湖店）'
$ws4.Range("E13").Value = '2024.11.30 10:00-11.30 21:00'
$ws4.Range("F13").Value = 98
$ws4.Range("G13").Value = 55
$ws4.Range("H13").Value = 'https://show.bilibili.com/platform/detail.html?id=93322'
$ws4.Range("I13").Value = '//i1.hdslb.com/bfs/openplatform/202409/Tu5YLbGx1727179854562.jpeg'
# Row 14
$ws4.Range("A14").Value = 13
$ws4.Range("B14").NumberFormat = "@"
$ws4.Range("B14").Value = '2024-12-07'
$ws4.Range("C14").Value = '合肥·一生必听的古典系列《钟》—超技钢琴曲炫彩音乐会'
$ws4.Range("D14").Value = '徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院'
$ws4.Range("E14").Value = '2024.12.07 19:30-12.07 21:00'
$ws4.Range("F14").Value = 2
$ws4.Range("G14").Value = 72
$ws4.Range("H14").Value = 'https://show.bilibili.com/platform/detail.html?id=91608'
$ws4.Range("I14").Value = '//i0.hdslb.com/bfs/openplatform/202408/wiLiWoeM1725005636569.jpeg'
# Row 15
$ws4.Range("A15").Value = 14
$ws4.Range("B15").NumberFormat = "@"
$ws4.Range("B15").Value = '2024-12-07'
$ws4.Range("C15").Value = '合肥·心动恋章·冬日序国乙&代号鸢同人only'
$ws4.Range("D15").Value = '上海路与迎淮路交口向东200米(云峯中心一楼) 费加罗宴会艺术中心(省府店)'
$ws4.Range("E15").Value = '2024.12.07 12:00-12.07 21:00'
$ws4.Range("F15").Value = 344
$ws4.Range("G15").Value = 77
$ws4.Range("H15").Value = 'https://show.bilibili.com/platform/detail.html?id=93319'
$ws4.Range("I15").Value = '//i0.hdslb.com/bfs/openplatform/202409/KtMLL8ZO1727684987784.jpeg'
# Row 16
$ws4.Range("A16").Value = 15
$ws4.Range("B16").NumberFormat = "@"
$ws4.Range("B16").Value = '2024-12-08'
$ws4.Range("C16").Value = '合肥·星光国潮动漫游戏嘉年华'
$ws4.Range("D16").Value = '北二环与新蚌埠路交汇处 蓝金湾大酒店'
$ws4.Range("E16").Value = '2024.12.08 10:00-12.08 17:00'
$ws4.Range("F16").Value = 8
$ws4.Range("G16").Value = 39.9
$ws4.Range("H16").Value = 'https://show.bilibili.com/platform/detail.html?id=93801'
$ws4.Range("I16").Value = '//i0.hdslb.com/bfs/openplatform/202410/ubX6VZ841729253636894.png'
# Row 17
$ws4.Range("A17").Value = 16
$ws4.Range("B17").NumberFormat = "@"
$ws4.Range("B17").Value = '2025-02-03'
$ws4.Range("C17").Value = '合肥·皖萌次元青年文化节'
$ws4.Range("D17").Value = '凤淮路与公园路交叉口南行50米路西 庐阳区全民健身中心'
$ws4.Range("E17").Value = '2025.02.03 10:00-02.04 17:30'
$ws4.Range("F17").Value = 63
$ws4.Range("G17").Value = 39.9
$ws4.Range("H17").Value = 'https://show.bilibili.com/platform/detail.html?id=93589'
$ws4.Range("I17").Value = '//i0.hdslb.com/bfs/openplatform/202410/GjWiXfOf1729133962063.jpeg'
